$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the stats for row 25 (2025-12) with newly reported values
$ws.Range("B25").Value = 6438
$ws.Range("D25").Value = 5996656

# Recompute derived metrics to match the updated raw values
$ws.Range("E25").Value = 931.4470332401366
$ws.Range("F25").Value = 9.285350534713977
$ws.Range("H25").Value = 25.58846924068337
